# Artillery scenario calcs: add an "ACTUAL VU created" column and a new
# rampup/"cloud batch file" scenario block (rows 13 & 15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column before N ("ACTUAL RPS max" etc. shift right) ---
$ws.Columns("N:N").Insert()

# Insert() stamps every row that had data in M with an empty, style-carried
# "N" cell. Only row 1 (header) and the new row 13 actually want content
# there, so drop the incidental empties it created in the pre-existing rows.
$ws.Range("N3").Clear()
$ws.Range("N4").Clear()
$ws.Range("N5").Clear()
$ws.Range("N7").Clear()
$ws.Range("N10").Clear()

# New header for the inserted column.
$ws.Range("N1").Value = "ACTUAL VU created"

# --- 2. New scenario block: rows 13 (data) and 15 (note), following the
#        same layout as the existing phase blocks above. ---
$ws.Rows("13").RowHeight = 30
$ws.Rows("15").RowHeight = 45

$ws.Range("B13").Value = "Rampup"

$ws.Range("C13").Value = 2
$ws.Range("D13").Value = 10
$ws.Range("E13").Value = 50
$ws.Range("F13").Formula = "=AVERAGE(C13,D13)*E13"
$ws.Range("G13").Value = 3
$ws.Range("H13").Value = 4
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 20
$ws.Range("K13").Formula = "=(G13*H13+I13)*J13"
$ws.Range("L13").Formula = "=(F13*G13*J13)/K13"
$ws.Range("M13").Formula = "=L13*10"

$ws.Range("N13").Value = 300
$ws.Range("N13").NumberFormat = "0.00"

$ws.Range("O13").Value = 750
$ws.Range("P13").Value = 0

# Shared-string insertion order follows the author's actual authoring
# sequence (note column first, then the run-id label, then the closing note).
$ws.Range("Q13").Value = "rampup as expected"
$ws.Range("A13").Value = "asciiArt_2024-05-23_run01"
$ws.Range("A15").Value = "PW (Playwright) run01"

# --- 3. Selection matches the author's final cursor position ---
$ws.Range("C15").Select()
